# Updates cryptos list cell values (prices and 1h volume %) for Thu Sep 12 17:28:17 UTC 2024 refresh.
# Values are written with a leading apostrophe (via string concatenation) so Excel stores them
# as literal text, matching the original inline-string cells (prevents things like "1.00" -> 1,
# "17.50" -> 17.5, or "58.128.33" misparsing).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '58.258.46'
$ws.Range('E2').Value = "'" + '  +1.90%  '
$ws.Range('D3').Value = "'" + '2.342.78'
$ws.Range('E3').Value = "'" + '  +0.62%  '
$ws.Range('E4').Value = "'" + '  -0.15%  '
$ws.Range('D5').Value = "'" + '542.65'
$ws.Range('E5').Value = "'" + '  +2.61%  '
$ws.Range('D6').Value = "'" + '134.91'
$ws.Range('E6').Value = "'" + '  +1.79%  '
$ws.Range('D7').Value = "'" + '0.999'
$ws.Range('E7').Value = "'" + '  +0.54%  '
$ws.Range('D8').Value = "'" + '0.562'
$ws.Range('E8').Value = "'" + '  +5.09%  '
$ws.Range('D9').Value = "'" + '0.102'
$ws.Range('E9').Value = "'" + '  +0.37%  '
$ws.Range('D10').Value = "'" + '5.59'
$ws.Range('E10').Value = "'" + '  +5.19%  '
$ws.Range('E11').Value = "'" + '  -0.55%  '
$ws.Range('D12').Value = "'" + '0.356'
$ws.Range('E12').Value = "'" + '  +3.21%  '
$ws.Range('D13').Value = "'" + '23.82'
$ws.Range('E13').Value = "'" + '  +1.08%  '
$ws.Range('D14').Value = "'" + '2.757.51'
$ws.Range('E14').Value = "'" + '  +0.29%  '
$ws.Range('D15').Value = "'" + '58.185.73'
$ws.Range('E15').Value = "'" + '  +1.72%  '
$ws.Range('E16').Value = "'" + '  +0.35%  '
$ws.Range('D17').Value = "'" + '2.335.46'
$ws.Range('E17').Value = "'" + '  -0.40%  '
$ws.Range('D18').Value = "'" + '10.69'
$ws.Range('E18').Value = "'" + '  +2.25%  '
$ws.Range('D19').Value = "'" + '333.29'
$ws.Range('E19').Value = "'" + '  -1.08%  '
$ws.Range('D20').Value = "'" + '4.25'
$ws.Range('E20').Value = "'" + '  +1.82%  '
$ws.Range('E21').Value = "'" + '  -3.26%  '
$ws.Range('E22').Value = "'" + '  +0.07%  '
$ws.Range('E23').Value = "'" + '  +0.12%  '
$ws.Range('D24').Value = "'" + '62.92'
$ws.Range('E24').Value = "'" + '  +1.94%  '
$ws.Range('E25').Value = "'" + '  +1.82%  '
$ws.Range('D26').Value = "'" + '8.53'
$ws.Range('E26').Value = "'" + '  -4.27%  '
$ws.Range('D27').Value = "'" + '0.999'
$ws.Range('E27').Value = "'" + '  +0.88%  '
$ws.Range('E28').Value = "'" + '  +6.05%  '
$ws.Range('E29').Value = "'" + '  +2.97%  '
$ws.Range('D30').Value = "'" + '170.87'
$ws.Range('E30').Value = "'" + '  +0.55%  '
$ws.Range('D31').Value = "'" + '0.0₃0736'
$ws.Range('E31').Value = "'" + '  +1.15%  '
$ws.Range('D32').Value = "'" + '6.11'
$ws.Range('B33').Value = "'" + 'SuiNetwork'
$ws.Range('C33').Value = "'" + 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D33').Value = "'" + '1.03'
$ws.Range('E33').Value = "'" + '  +13.19%  '
$ws.Range('B34').Value = "'" + 'EthereumClassic'
$ws.Range('C34').Value = "'" + 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = "'" + '18.45'
$ws.Range('E34').Value = "'" + '  -0.54%  '
$ws.Range('E35').Value = "'" + '  +0.04%  '
$ws.Range('D36').Value = "'" + '4.26'
$ws.Range('E36').Value = "'" + '  +5.82%  '
$ws.Range('D37').Value = "'" + '1.00'
$ws.Range('E37').Value = "'" + '  +0.89%  '
$ws.Range('E38').Value = "'" + '  -1.79%  '
$ws.Range('D39').Value = "'" + '1.65'
$ws.Range('E39').Value = "'" + '  +3.34%  '
$ws.Range('D40').Value = "'" + '39.08'
$ws.Range('E40').Value = "'" + '  +0.31%  '
$ws.Range('D41').Value = "'" + '143.54'
$ws.Range('E41').Value = "'" + '  -3.30%  '
$ws.Range('B42').Value = "'" + 'PolygonEcosystemToken'
$ws.Range('C42').Value = "'" + 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').Value = "'" + '0.376'
$ws.Range('E42').Value = "'" + '  -0.60%  '
$ws.Range('B43').Value = "'" + 'Filecoin'
$ws.Range('C43').Value = "'" + 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = "'" + '3.64'
$ws.Range('E43').Value = "'" + '  +1.17%  '
$ws.Range('D44').Value = "'" + '286.26'
$ws.Range('E44').Value = "'" + '  -0.13%  '
$ws.Range('D45').Value = "'" + '0.0940'
$ws.Range('E45').Value = "'" + '  +0.91%  '
$ws.Range('E46').Value = "'" + '  +2.33%  '
$ws.Range('D47').Value = "'" + '0.0503'
$ws.Range('E47').Value = "'" + '  +0.02%  '
$ws.Range('D48').Value = "'" + '0.564'
$ws.Range('E48').Value = "'" + '  +0.56%  '
$ws.Range('B49').Value = "'" + 'Polygon'
$ws.Range('C49').Value = "'" + 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D49').Value = "'" + '0.385'
$ws.Range('E49').Value = "'" + '  +1.74%  '
$ws.Range('B50').Value = "'" + 'VeChain'
$ws.Range('C50').Value = "'" + 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = "'" + '0.0218'
$ws.Range('E50').Value = "'" + '  +0.29%  '
$ws.Range('B51').Value = "'" + 'EnergySwap'
$ws.Range('C51').Value = "'" + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'" + '17.50'
$ws.Range('E51').Value = "'" + '  +0.78%  '
